$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos table (GitHub Actions scheduled update): new Price (D)
# and Volume(1h) (E) readings for most rows, plus two rows (34/35) where
# ImmutableX/USDe swapped rank order, so their Coin/Link/Price/Volume cells
# are all replaced.
#
# Column D holds prices as plain text (e.g. "59.136.07", or values with a
# significant trailing zero like "6.70"/"10.90"). Where the new price string
# would otherwise be auto-parsed as a number by Excel, it is entered with a
# leading apostrophe (exactly like typing '6.70 into a cell) and the cell's
# style is reset to Normal afterwards so no stray NumberFormat/quote-prefix
# style is left on the cell.
$ws.Range('D2').Value = '59.136.07'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '2.521.71'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'536.35"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = "'137.85"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = "'0.568"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = '2.518.09'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').Value = "'5.36"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = "'0.349"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').Value = '2.968.58'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = "'23.01"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('D16').Value = '59.057.86'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '2.518.22'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = "'11.11"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').Value = "'325.86"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = "'5.96"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.43%  '
$ws.Range('D24').Value = "'65.83"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = "'6.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  +6.18%  '
$ws.Range('D33').Value = "'163.77"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'1.47"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').Value = "'18.48"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = "'36.59"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').Value = "'3.64"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.47%  '
$ws.Range('D42').Value = "'286.31"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('D43').Value = "'5.21"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = "'132.67"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.66%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').Value = "'0.605"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').Value = "'10.90"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = "'0.0511"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('E51').Value = '  -2.39%  '
